# Ships characteristics workbook update
# - adds a third ("My version") table to sheet 1 (rows 18-27)
# - tweaks a couple of column widths
# - updates the selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# ---------------------------------------------------------------------------
# Column width tweaks (column B gets narrower, new column G gets a width)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 8.43
$ws.Columns.Item(7).ColumnWidth = 16.66

# ---------------------------------------------------------------------------
# New table - written in a deliberate order so that new shared strings end
# up in the same order as the authored workbook:
#   17 = "Cutter", 18 = "My version", 19 = "Paddle speed", 20 = "Jolly boat"
# ---------------------------------------------------------------------------

# "Cutter" must become shared-string #17 -> write it first
$ws.Range("A22").Value = "Cutter"

# Section title -> shared-string #18
$ws.Range("A18").Value = "My version"

# Header row (row 19) -> last new header label "Paddle speed" is #19
$ws.Range("A19").Value = "Model"
$ws.Range("B19").Value = "Model number"
$ws.Range("C19").Value = "Price"
$ws.Range("D19").Value = "Health"
$ws.Range("E19").Value = "Speed"
$ws.Range("F19").Value = "Speed towards wind"
$ws.Range("G19").Value = "Paddle speed"
$ws.Range("H19").Value = "Crew"
$ws.Range("I19").Value = "Cannons (initial)"
$ws.Range("J19").Value = "Cannons (max)"
$ws.Range("K19").Value = "Cargo space"

# Row 20 - "Jolly boat" -> #20
$ws.Range("A20").Value = "Jolly boat"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 500
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 15
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 1
$ws.Range("K20").Value = 1

# Row 21 - Sloop
$ws.Range("A21").Value = "Sloop"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 18000
$ws.Range("D21").Value = 50
$ws.Range("E21").NumberFormat = "0.00"
$ws.Range("E21").Value = 4.05
$ws.Range("F21").Value = 2.355
$ws.Range("G21").Value = 0.9
$ws.Range("H21").Value = 66
$ws.Range("I21").Value = 12
$ws.Range("J21").Value = 16
$ws.Range("K21").Value = 2

# Row 22 - Cutter (label already written above, fill the rest of the row)
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = 25000
$ws.Range("D22").Value = 80
$ws.Range("E22").Value = 3.7
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 16
$ws.Range("J22").Value = 20
$ws.Range("K22").Value = 2

# Row 23 - Caravel
$ws.Range("A23").Value = "Caravel"
$ws.Range("B23").Value = 3
$ws.Range("C23").Value = 38500
$ws.Range("D23").Value = 115
$ws.Range("E23").Value = 2.85
$ws.Range("F23").Value = 1.05
$ws.Range("G23").Value = 0.3
$ws.Range("H23").Value = 175
$ws.Range("I23").Value = 24
$ws.Range("J23").Value = 30
$ws.Range("K23").Value = 4

# Row 24 - Brig
$ws.Range("A24").Value = "Brig"
$ws.Range("B24").Value = 4
$ws.Range("C24").Value = 41500
$ws.Range("D24").Value = 135
$ws.Range("E24").Value = 4.35
$ws.Range("F24").Value = 1.05
$ws.Range("G24").Value = 0.4
$ws.Range("H24").Value = 177
$ws.Range("I24").Value = 16
$ws.Range("J24").Value = 24
$ws.Range("K24").Value = 3

# Row 25 - Galeon
$ws.Range("A25").Value = "Galeon"
$ws.Range("B25").Value = 5
$ws.Range("C25").Value = 100000
$ws.Range("D25").Value = 280
$ws.Range("E25").Value = 2.25
$ws.Range("F25").Value = 0.78
$ws.Range("G25").Value = 0.2
$ws.Range("H25").Value = 448
$ws.Range("I25").Value = 32
$ws.Range("J25").Value = 36
$ws.Range("K25").Value = 7

# Row 26 - Fregat
$ws.Range("A26").Value = "Fregat"
$ws.Range("B26").Value = 6
$ws.Range("C26").Value = 150000
$ws.Range("D26").Value = 250
$ws.Range("E26").Value = 4.65
$ws.Range("F26").Value = 1.17
$ws.Range("G26").Value = 0.5
$ws.Range("H26").Value = 323
$ws.Range("I26").Value = 32
$ws.Range("J26").Value = 46
$ws.Range("K26").Value = 4

# Row 27 - Battleship
$ws.Range("A27").Value = "Battleship"
$ws.Range("B27").Value = 7
$ws.Range("C27").Value = 250000
$ws.Range("D27").Value = 380
$ws.Range("E27").Value = 3.75
$ws.Range("F27").Value = 0.42
$ws.Range("G27").Value = 0.3
$ws.Range("H27").Value = 571
$ws.Range("I27").Value = 32
$ws.Range("J27").Value = 66
$ws.Range("K27").Value = 5

# ---------------------------------------------------------------------------
# View state: scroll so row 4 is at the top and select G26
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 2
$ws.Range("G26").Select()
